$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = $null
$ws.Range("N6").Value = $null

$ws.Range("H116").Value = 6030409.5
$ws.Range("I116").Value = 6298301
$ws.Range("J116").Value = 2850
$ws.Range("K116").Value = 6298301
$ws.Range("L116").Value = 2850
$ws.Range("M116").Value = -6294859
$ws.Range("N116").Value = -9734

$ws.Range("H132").Value = 2480.3252
$ws.Range("I132").Value = 2129.621
$ws.Range("J132").Value = 3841.8823
$ws.Range("K132").Value = 6388.863
$ws.Range("L132").Value = 11525.6469
$ws.Range("M132").Value = -3858.863
$ws.Range("N132").Value = -16585.6469

$ws.Range("H137").Value = 1827.5883
$ws.Range("I137").Value = 1428.7333
$ws.Range("J137").Value = 2142.4736
$ws.Range("K137").Value = 4286.199900000001
$ws.Range("L137").Value = 6427.4208
$ws.Range("M137").Value = -1736.199900000001
$ws.Range("N137").Value = -11527.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1182.0646
$ws.Range("I61").Value = 1106.4783
$ws.Range("J61").Value = 1399.375
$ws.Range("K61").Value = 1106.4783
$ws.Range("L61").Value = 1399.375
$ws.Range("M61").Value = -894.4783
$ws.Range("N61").Value = -1823.375

$ws.Range("H74").Value = 7578552.5
$ws.Range("I74").Value = 10872205
$ws.Range("J74").Value = 3152.6
$ws.Range("K74").Value = 10872205
$ws.Range("L74").Value = 3152.6
$ws.Range("M74").Value = -10871331
$ws.Range("N74").Value = -4900.6

$ws.Range("H77").Value = 7578552.5
$ws.Range("I77").Value = 10872205
$ws.Range("J77").Value = 3152.6
$ws.Range("K77").Value = 54361025
$ws.Range("L77").Value = 15763
$ws.Range("M77").Value = -54356657
$ws.Range("N77").Value = -24499

$ws.Range("H132").Value = 1515.4108
$ws.Range("I132").Value = 1252.3617
$ws.Range("J132").Value = 2889.111
$ws.Range("K132").Value = 3757.0851
$ws.Range("L132").Value = 8667.332999999999
$ws.Range("M132").Value = -1227.0851
$ws.Range("N132").Value = -13727.333

$ws.Range("H136").Value = 1182.0646
$ws.Range("I136").Value = 1106.4783
$ws.Range("J136").Value = 1399.375
$ws.Range("K136").Value = 3319.4349
$ws.Range("L136").Value = 4198.125
$ws.Range("M136").Value = -769.4349000000002
$ws.Range("N136").Value = -9298.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10832

$ws.Range("H76").Value = 9157
$ws.Range("J76").Value = 9157
$ws.Range("L76").Value = 9157
$ws.Range("N76").Value = -9787

$ws.Range("H79").Value = 9157
$ws.Range("J79").Value = 9157
$ws.Range("L79").Value = 9157
$ws.Range("N79").Value = -11341

$ws.Range("H88").Value = 42750
$ws.Range("J88").Value = 42750
$ws.Range("L88").Value = 42750
$ws.Range("N88").Value = -43562

$ws.Range("H91").Value = 42750
$ws.Range("J91").Value = 42750
$ws.Range("L91").Value = 42750
$ws.Range("N91").Value = -45558

$ws.Range("H132").Value = 26993.334
$ws.Range("J132").Value = 26993.334
$ws.Range("L132").Value = 26993.334
$ws.Range("N132").Value = -37113.334

$ws.Range("H134").Value = 2156.8845
$ws.Range("I134").Value = 1809.9445
$ws.Range("J134").Value = 2937.5
$ws.Range("K134").Value = 5429.833500000001
$ws.Range("L134").Value = 8812.5
$ws.Range("M134").Value = -2894.833500000001
$ws.Range("N134").Value = -13882.5

$ws.Range("H135").Value = 28685.715
$ws.Range("J135").Value = 28685.715
$ws.Range("L135").Value = 28685.715
$ws.Range("N135").Value = -38825.715

$ws.Range("H137").Value = 46933.332
$ws.Range("J137").Value = 46933.332
$ws.Range("L137").Value = 46933.332
$ws.Range("N137").Value = -57133.332

$ws.Range("H138").Value = 47220
$ws.Range("J138").Value = 47220
$ws.Range("L138").Value = 47220
$ws.Range("N138").Value = -57500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 19969
$ws.Range("J74").Value = 19969
$ws.Range("L74").Value = 19969
$ws.Range("N74").Value = -21717

$ws.Range("H77").Value = 19969
$ws.Range("J77").Value = 19969
$ws.Range("L77").Value = 59907
$ws.Range("N77").Value = -68643

$ws.Range("H134").Value = 2577.3713
$ws.Range("I134").Value = 3044.1738
$ws.Range("J134").Value = 1682.6666
$ws.Range("K134").Value = 9132.5214
$ws.Range("L134").Value = 5047.9998
$ws.Range("M134").Value = -6597.5214
$ws.Range("N134").Value = -10117.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 428.90475
$ws.Range("I5").Value = 358.51352
$ws.Range("J5").Value = 949.8
$ws.Range("K5").Value = 1075.54056
$ws.Range("L5").Value = 2849.4
$ws.Range("M5").Value = -963.5405600000001
$ws.Range("N5").Value = -3073.4

$ws.Range("H92").Value = 711.4286
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 836
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 2508
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -5004

$ws.Range("H135").Value = 428.90475
$ws.Range("I135").Value = 358.51352
$ws.Range("J135").Value = 949.8
$ws.Range("K135").Value = 3226.62168
$ws.Range("L135").Value = 8548.199999999999
$ws.Range("M135").Value = -691.6216800000002
$ws.Range("N135").Value = -13618.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4196.615
$ws.Range("I122").Value = 3631.6875
$ws.Range("J122").Value = 5100.5
$ws.Range("K122").Value = 10895.0625
$ws.Range("L122").Value = 15301.5
$ws.Range("M122").Value = -8445.0625
$ws.Range("N122").Value = -20201.5

$ws.Range("H132").Value = 3031.7222
$ws.Range("I132").Value = 1997
$ws.Range("J132").Value = 4066.4443
$ws.Range("K132").Value = 5991
$ws.Range("L132").Value = 12199.3329
$ws.Range("M132").Value = -3461
$ws.Range("N132").Value = -17259.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 17369394
$ws.Range("I132").Value = 44659756
$ws.Range("J132").Value = 2797.8635
$ws.Range("K132").Value = 133979268
$ws.Range("L132").Value = 8393.5905
$ws.Range("M132").Value = -133976738
$ws.Range("N132").Value = -13453.5905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 38463310
$ws.Range("I81").Value = 71430140
$ws.Range("J81").Value = 2008.5
$ws.Range("K81").Value = 142860280
$ws.Range("L81").Value = 4017
$ws.Range("M81").Value = -142859219
$ws.Range("N81").Value = -6139

$ws.Range("H84").Value = 38463310
$ws.Range("I84").Value = 71430140
$ws.Range("J84").Value = 2008.5
$ws.Range("K84").Value = 714301400
$ws.Range("L84").Value = 20085
$ws.Range("M84").Value = -714296096
$ws.Range("N84").Value = -30693

$ws.Range("H132").Value = 1052.3673
$ws.Range("I132").Value = 564.32355
$ws.Range("K132").Value = 1692.97065
$ws.Range("M132").Value = 837.0293500000002
